$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: conversion rates text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.79 = 41351.62 pesos`n✅ 41351.62 pesos = 9.73 = 917.78 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet: updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 102.1
$ws2.Range("O10").Value = 4222
$ws2.Range("N12").Value = 4249.7
$ws2.Range("O12").Value = 94.32
